$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 362, shifting the existing rows 362:393 down to 363:394.
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row with the new weekly price-report record.
$ws.Cells.Item(362, 1).Value = 9
$ws.Cells.Item(362, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(362, 3).Value = "Metropolitana"
$ws.Cells.Item(362, 4).Value = 45194
$ws.Cells.Item(362, 5).Value = 13
$ws.Cells.Item(362, 6).Value = 100112026
$ws.Cells.Item(362, 7).Value = "Haba"
$ws.Cells.Item(362, 8).Value = "Sin especificar"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 52
$ws.Cells.Item(362, 11).Value = 12000
$ws.Cells.Item(362, 12).Value = 13000
$ws.Cells.Item(362, 13).Value = 12500
$ws.Cells.Item(362, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(362, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(362, 16).Value = 500
$ws.Cells.Item(362, 17).Value = 25
$ws.Cells.Item(362, 18).Value = "Hortaliza"
